$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/crossover-indicator"
$ws1.Range("B3").Value = "8.0.0"
$ws1.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$ws1.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$ws2 = $wb.Worksheets.Item("Elements")

# Row 2 (Extension): clear the Constraint(s) cell (AI2) - it was misplaced here
$ws2.Range("AI2").Value = ""

# Row 4 (Extension.extension): populate slicing info + constraint(s) that belong here
$ws2.Range("AA4").Value = "value:url}`n"
$ws2.Range("AB4").Value = "Extensions are always sliced by (at least) url"
$ws2.Range("AD4").Value = "open"
$ws2.Range("AI4").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 5 (Extension.url): Fixed Value reflects the updated canonical extension URL
$ws2.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/crossover-indicator"

# Row 6 (Extension.value[x]): Base Max cell value
$ws2.Range("AG6").Value = "string`n"
